$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

# Update Runmode (column E) for rows 9 through 21 from "No" to "Yes"
$ws.Range("E9:E21").Value = "Yes"

# Update the selection to match the new state (E8:E21)
$ws.Range("E8:E21").Select()
